# Insert a new weekly record row into the "Poroto verde" price sheet.
# This pushes the existing rows 520..547 down to 521..548 and populates
# the newly opened row 520 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 520 (shifts 520:547 -> 521:548).
$ws.Rows.Item(520).Insert()

# Fill the new row 520 with the new weekly observation.
$ws.Cells.Item(520, 1).Value = 9
$ws.Cells.Item(520, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(520, 3).Value = "Metropolitana"
$ws.Cells.Item(520, 4).Value = 44753
$ws.Cells.Item(520, 5).Value = 13
$ws.Cells.Item(520, 6).Value = 100112031
$ws.Cells.Item(520, 7).Value = "Poroto verde"
$ws.Cells.Item(520, 8).Value = "Magnum"
$ws.Cells.Item(520, 9).Value = "Primera"
$ws.Cells.Item(520, 10).Value = 52
$ws.Cells.Item(520, 11).Value = 34000
$ws.Cells.Item(520, 12).Value = 35000
$ws.Cells.Item(520, 13).Value = 34500
$ws.Cells.Item(520, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(520, 15).Value = "Perú"
$ws.Cells.Item(520, 16).Value = 1380
$ws.Cells.Item(520, 17).Value = 25
$ws.Cells.Item(520, 18).Value = "Hortaliza"
